$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The Left/Top/Width/Height properties on the PowerPoint shape object model
# are expressed in points; the target geometry below comes from the OOXML
# diff, which is in EMU (1 pt = 12700 EMU).
$emuPerPt = 12700

# A tiny positive nudge (well under a twentieth of an EMU once converted
# back) that keeps the single-precision round-trip PowerPoint's shape
# geometry uses from truncating down to one EMU less than the target.
$fudge = 0.00001

# --- Resize/reposition the existing Fresnel-equation picture -------------
$pic = $s.Shapes.Item("Picture 18")
$pic.Left = 4214297 / $emuPerPt
$pic.Top = 4825033 / $emuPerPt
$pic.Width = (1936453 / $emuPerPt) + $fudge
$pic.Height = 1542456 / $emuPerPt

# --- Add the new caption textbox ------------------------------------------
# Burn through the next few auto-assigned shape-id slots (mirroring ids
# PowerPoint already spent elsewhere in this deck) so the textbox we keep
# lands on id 21 / "TextBox 20", matching the original authoring session.
for ($i = 0; $i -lt 9; $i++) {
    $dummy = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $dummy.Delete()
}

# Duplicate a sibling caption textbox rather than building one from scratch
# so the new shape naturally inherits the same body/paragraph formatting
# (wrap="square" rtlCol="0", spAutoFit, lstStyle, dirty="0" run props, ...).
$src = $s.Shapes.Item("TextBox 16")
$tb = $src.Duplicate()
$tb.Name = "TextBox 20"
$tb.TextFrame.TextRange.Text = "Applies to any linear, isotropic, and homogeneous medium."
$tb.Left = 6501934 / $emuPerPt
$tb.Top = 5813742 / $emuPerPt
$tb.Width = 4058060 / $emuPerPt
$tb.Height = (646331 / $emuPerPt) + $fudge
